# Applies the "Fix 15" commit: adds a new FHIR extension slice
# "Medication.ingredient.extension:IsVehicle" to the Elements sheet of the
# fr-medication-compound2 StructureDefinition workbook, turns the parent
# "Medication.ingredient.extension" row into a sliced-extension header row,
# and bumps the Metadata "Date" property.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump the publication Date property (B8)
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-01-11T17:10:04+00:00"

# ---------------------------------------------------------------------
# 2. Elements sheet: insert the new slice row and update the parent row
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Insert a new row at 24 (pushes the old rows 24-33 down to 25-34) and
# copy formatting + values from row 23 (the extension parent row) as a
# starting point so the new row inherits the same cell style (s="2").
$ws.Rows.Item(24).Insert()
$ws.Range("A23:AN23").Copy($ws.Range("A24:AN24"))

# --- Row 23 ("Medication.ingredient.extension") becomes the generic
#     sliced-extension header, mirroring how "Medication.code.extension"
#     (row 14) is already modelled.
$ws.Range("D23").Value = ""
$ws.Range("L23").Value = "Extension"
$ws.Range("M23").Value = "An Extension"
$ws.Range("N23").Value = ""
$ws.Range("AB23").Value = "value:url}" + [char]10
$ws.Range("AC23").Value = ""
$ws.Range("AE23").Value = "open"
$ws.Range("AL23").Value = ""

# --- Row 24 (new): "Medication.ingredient.extension:IsVehicle" slice
$ws.Range("A24").Value = "Medication.ingredient.extension:IsVehicle"
$ws.Range("B24").Value = "Medication.ingredient.extension"
$ws.Range("C24").Value = "IsVehicle"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = "0"
$ws.Range("G24").Value = "1"
$ws.Range("H24").Value = ""
$ws.Range("I24").Value = ""
$ws.Range("J24").Value = ""
$ws.Range("K24").Value = "Extension {https://hl7.fr/fhir/fr/medication/StructureDefinition/fr-is-vehicle}" + [char]10
$ws.Range("L24").Value = "Medication component which is the vehicle of the compound Medication"
$ws.Range("M24").Value = "Medication component which is the vehicle of the compound Medication"
$ws.Range("N24").Value = ""
$ws.Range("O24").Value = ""
$ws.Range("P24").Value = ""
$ws.Range("Q24").Value = ""
$ws.Range("R24").Value = ""
$ws.Range("S24").Value = ""
$ws.Range("T24").Value = ""
$ws.Range("U24").Value = ""
$ws.Range("V24").Value = ""
$ws.Range("W24").Value = ""
$ws.Range("X24").Value = ""
$ws.Range("Y24").Value = ""
$ws.Range("Z24").Value = ""
$ws.Range("AA24").Value = ""
$ws.Range("AB24").Value = ""
$ws.Range("AC24").Value = ""
$ws.Range("AD24").Value = ""
$ws.Range("AE24").Value = ""
$ws.Range("AF24").Value = "Element.extension"
$ws.Range("AG24").Value = "0"
$ws.Range("AH24").Value = "*"
$ws.Range("AI24").Value = "ele-1" + [char]10
$ws.Range("AJ24").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$ws.Range("AK24").Value = ""
$ws.Range("AL24").Value = ""
$ws.Range("AM24").Value = ""
$ws.Range("AN24").Value = ""

# ---------------------------------------------------------------------
# 3. Column A is now a little wider because of the longer ID string
#    "Medication.ingredient.extension:IsVehicle" - widen it to the
#    book's recorded best-fit width.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 38.6

Write-Output "Applied IsVehicle extension slice edit"
